# Fixing issues for TestCase_E5 and TestCase_E6
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Rows 2-5 (TestCase_E1..E4): Results column SKIP -> PASS
$ws.Range("E2:E5").Value = "PASS"

# Row 6 (TestCase_E5): Jira id TBD-01 -> OPQA-262
$ws.Range("B6").Value = "OPQA-262"

# Row 7 (TestCase_E6): Jira id TBD-02 -> OPQA-264
$ws.Range("B7").Value = "OPQA-264"

# Update the view: scroll so column D is the left-most visible column,
# and select F8 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F8").Select()
